# Update building block types
#  - bump template Version 1.0.2 -> 1.0.3
#  - rename "Parameter [library preparation kit]" -> "Component [library preparation kit]"
#  - rename "Parameter [next generation sequencing instrument model]" -> "Component [next generation sequencing instrument model]"
#  - swap the ontology term source/accession for the instrument model row from
#    OBI / https://bioregistry.io/OBI:0002049 to EFO / https://bioregistry.io/EFO:0008563

$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump version ---
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.3"

# --- next_generation_sequencing sheet: rename building block headers ---
$wsData = $wb.Worksheets.Item("next_generation_sequencing")
$wsData.Range("K1").Value = "Component [library preparation kit]"
$wsData.Range("N1").Value = "Component [next generation sequencing instrument model]"

# --- update the instrument-model ontology reference on the data row ---
$wsData.Range("O2").Value = "EFO"
$wsData.Range("P2").Value = "https://bioregistry.io/EFO:0008563"
